# Snap Share implementation across TP, PP, SC
# Remove rookies who are no longer tracked, and normalize a few player
# names by dropping their suffix (III / Jr.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Players to drop from the list entirely (entire row removed).
$playersToRemove = @(
    "Donovan Peoples-Jones",
    "Laviska Shenault",
    "Trey Sermon",
    "Terrace Marshall",
    "Amari Rodgers",
    "Jelani Woods",
    "Calvin Austin III"
)

# Players whose display name should be normalized (suffix dropped).
# Processed in this specific order (John Metchie, Brian Robinson, then
# Kenneth Walker) so new entries land in the same order in the workbook.
$renameOrder = @(
    "John Metchie III",
    "Brian Robinson Jr.",
    "Kenneth Walker III"
)
$renameMap = @{
    "Kenneth Walker III"  = "Kenneth Walker"
    "John Metchie III"    = "John Metchie"
    "Brian Robinson Jr."  = "Brian Robinson"
}

$lastRow = $ws.UsedRange.Rows.Count

# First, rename the players that need their suffix dropped.
foreach ($oldName in $renameOrder) {
    for ($r = 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $val = $cell.Text
        if ($val -eq $oldName) {
            $cell.Value = $renameMap[$oldName]
        }
    }
}

# Then delete the rows for removed players, from the bottom up so row
# numbers of not-yet-processed rows stay stable.
for ($r = $lastRow; $r -ge 1; $r--) {
    $val = $ws.Cells.Item($r, 1).Text
    if ($playersToRemove -contains $val) {
        $ws.Rows.Item($r).Delete()
    }
}

# Update the view to reflect the new, shorter list.
$ws.Range("A50").Select()
